$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 124
$ws1.Range("F3").Value = 16

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 124
$ws4.Range("F3").Value = 16
